# Rename the worksheet "Data_Final" to "Data-Final" to avoid loading
# errors when the file is opened from Jupyter.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data_Final")
$ws.Name = "Data-Final"
